$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) cells that are being updated, so that
# numeric-looking strings (e.g. "0.9980", "22.90") are preserved exactly as
# text instead of being parsed/rounded into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Update Price (D) values
$ws.Range("D2").Value = '29.073.21'
$ws.Range("D3").Value = '1.820.54'
$ws.Range("D4").Value = '0.9980'
$ws.Range("D5").Value = '241.31'
$ws.Range("D6").Value = '0.6147'
$ws.Range("D7").Value = '0.9995'
$ws.Range("D8").Value = '0.07326'
$ws.Range("D9").Value = '0.2884'
$ws.Range("D10").Value = '22.90'
$ws.Range("D11").Value = '0.07652'
$ws.Range("D12").Value = '1.811.81'
$ws.Range("D13").Value = '4.939'
$ws.Range("D14").Value = '0.6586'
$ws.Range("D15").Value = '81.62'
$ws.Range("D16").Value = '0.000008996'
$ws.Range("D17").Value = '5.821'
$ws.Range("D18").Value = '29.050.16'
$ws.Range("D19").Value = '2.066.94'
$ws.Range("D20").Value = '236.88'
$ws.Range("D22").Value = '0.9994'
$ws.Range("D23").Value = '7.101'
$ws.Range("D24").Value = '0.9997'
$ws.Range("D25").Value = '157.24'
$ws.Range("D26").Value = '0.1398'
$ws.Range("D27").Value = '8.415'
$ws.Range("D28").Value = '17.55'
$ws.Range("D29").Value = '1.483'
$ws.Range("D30").Value = '0.05553'
$ws.Range("D31").Value = '4.091'
$ws.Range("D32").Value = '4.085'
$ws.Range("D34").Value = '0.7330'
$ws.Range("D35").Value = '1.810'
$ws.Range("D36").Value = '1.129'
$ws.Range("D37").Value = '2.608'
$ws.Range("D38").Value = '2.824'
$ws.Range("D39").Value = '1.206.90'
$ws.Range("D41").Value = '6.351'
$ws.Range("D42").Value = '0.8915'
$ws.Range("D43").Value = '0.9998'
$ws.Range("D44").Value = '100.95'
$ws.Range("D45").Value = '1.969.69'
$ws.Range("D46").Value = '64.36'
$ws.Range("D47").Value = '0.5077'
$ws.Range("D49").Value = '0.3991'
$ws.Range("D50").Value = '9.028'

# Update Volume(1h) (E) values
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("E6").Value = '  -2.13%  '
$ws.Range("E7").Value = '  -0.31%  '
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("E16").Value = '  -3.97%  '
$ws.Range("E17").Value = '  -2.67%  '
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("E20").Value = '  +6.17%  '
$ws.Range("E21").Value = '  -1.32%  '
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("E25").Value = '  -1.82%  '
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("E27").Value = '  -0.87%  '
$ws.Range("E28").Value = '  -1.85%  '
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("E30").Value = '  -1.69%  '
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("E32").Value = '  -1.63%  '
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("E34").Value = '  -1.18%  '
$ws.Range("E35").Value = '  -1.45%  '
$ws.Range("E36").Value = '  -1.01%  '
$ws.Range("E37").Value = '  -2.42%  '
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("E39").Value = '  -1.02%  '
$ws.Range("E40").Value = '  -1.53%  '
$ws.Range("E41").Value = '  -2.81%  '
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("E43").Value = '  -0.25%  '
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("E46").Value = '  -2.00%  '
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("E48").Value = '  -4.81%  '
$ws.Range("E49").Value = '  -1.96%  '
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("E51").Value = '  -1.30%  '
